$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 799
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H43").Value = 931.6667
$ws.Range("J43").Value = 931.6667
$ws.Range("L43").Value = 931.6667
$ws.Range("N43").Value = -1069.6667

$ws.Range("H140").Value = 38926.668
$ws.Range("J140").Value = 38926.668
$ws.Range("L140").Value = 38926.668
$ws.Range("N140").Value = -49286.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2262.6
$ws.Range("I45").Value = 1320.5555
$ws.Range("J45").Value = 3675.6667
$ws.Range("K45").Value = 1320.5555
$ws.Range("L45").Value = 3675.6667
$ws.Range("M45").Value = -943.5554999999999
$ws.Range("N45").Value = -4429.6667

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040

$ws.Range("H132").Value = 2699.3333
$ws.Range("I132").Value = 2699.3333
$ws.Range("K132").Value = 8097.999899999999
$ws.Range("M132").Value = -5567.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4480.7144
$ws.Range("J86").Value = 6957.5
$ws.Range("L86").Value = 6957.5
$ws.Range("N86").Value = -9203.5

$ws.Range("H89").Value = 4480.7144
$ws.Range("J89").Value = 6957.5
$ws.Range("L89").Value = 34787.5
$ws.Range("N89").Value = -46019.5

$ws.Range("H94").Value = 598.6667
$ws.Range("I94").Value = 518.4
$ws.Range("K94").Value = 518.4
$ws.Range("M94").Value = -67.39999999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6922.4
$ws.Range("I7").Value = 9254.182000000001
$ws.Range("J7").Value = 510
$ws.Range("K7").Value = 9254.182000000001
$ws.Range("L7").Value = 510
$ws.Range("M7").Value = -9141.182000000001
$ws.Range("N7").Value = -736

$ws.Range("H21").Value = 31184.285
$ws.Range("J21").Value = 20723.889
$ws.Range("L21").Value = 20723.889
$ws.Range("N21").Value = -21193.889

$ws.Range("H29").Value = 3100
$ws.Range("J29").Value = 3100
$ws.Range("L29").Value = 3100
$ws.Range("N29").Value = -3686

$ws.Range("H31").Value = 9305.5
$ws.Range("I31").Value = 6200
$ws.Range("J31").Value = 9749.143
$ws.Range("K31").Value = 6200
$ws.Range("L31").Value = 9749.143
$ws.Range("M31").Value = -5905
$ws.Range("N31").Value = -10339.143

$ws.Range("H32").Value = 1316.6666
$ws.Range("I32").Value = 1316.6666
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1316.6666
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1000.6666
$ws.Range("N32").ClearContents()

$ws.Range("H33").Value = 696.7143
$ws.Range("I33").Value = 687.8333
$ws.Range("J33").Value = 750
$ws.Range("K33").Value = 687.8333
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -308.8333
$ws.Range("N33").Value = -1508

$ws.Range("H34").Value = 9305.5
$ws.Range("I34").Value = 6200
$ws.Range("J34").Value = 9749.143
$ws.Range("K34").Value = 6200
$ws.Range("L34").Value = 9749.143
$ws.Range("M34").Value = -5998
$ws.Range("N34").Value = -10153.143

$ws.Range("H38").Value = 6432.8
$ws.Range("I38").Value = 7897.25
$ws.Range("J38").Value = 575
$ws.Range("K38").Value = 7897.25
$ws.Range("L38").Value = 575
$ws.Range("M38").Value = -7520.25
$ws.Range("N38").Value = -1329

$ws.Range("H39").Value = 5460.2856
$ws.Range("I39").Value = 950.3333
$ws.Range("J39").Value = 8842.75
$ws.Range("K39").Value = 950.3333
$ws.Range("L39").Value = 8842.75
$ws.Range("M39").Value = -559.3333
$ws.Range("N39").Value = -9624.75

$ws.Range("H42").Value = 11638.75
$ws.Range("I42").Value = 11638.75
$ws.Range("K42").Value = 11638.75
$ws.Range("M42").Value = -11045.75

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H46").Value = 6432.8
$ws.Range("I46").Value = 7897.25
$ws.Range("J46").Value = 575
$ws.Range("K46").Value = 7897.25
$ws.Range("L46").Value = 575
$ws.Range("M46").Value = -7686.25
$ws.Range("N46").Value = -997

$ws.Range("H49").Value = 5460.2856
$ws.Range("I49").Value = 950.3333
$ws.Range("J49").Value = 8842.75
$ws.Range("K49").Value = 950.3333
$ws.Range("L49").Value = 8842.75
$ws.Range("M49").Value = -768.3333
$ws.Range("N49").Value = -9206.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 401
$ws.Range("I32").Value = 401
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1203
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -920
$ws.Range("N32").ClearContents()

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4500
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -4409
$ws.Range("N46").ClearContents()

$ws.Range("H59").Value = 595
$ws.Range("I59").Value = 595
$ws.Range("K59").Value = 1785
$ws.Range("M59").Value = -1245

$ws.Range("H68").Value = 2525.3333
$ws.Range("I68").Value = 1447.5
$ws.Range("J68").Value = 3064.25
$ws.Range("K68").Value = 4342.5
$ws.Range("L68").Value = 9192.75
$ws.Range("M68").Value = -3531.5
$ws.Range("N68").Value = -10814.75

$ws.Range("H71").Value = 2525.3333
$ws.Range("I71").Value = 1447.5
$ws.Range("J71").Value = 3064.25
$ws.Range("K71").Value = 13027.5
$ws.Range("L71").Value = 27578.25
$ws.Range("M71").Value = -8971.5
$ws.Range("N71").Value = -35690.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2800
$ws.Range("I22").Value = 2100
$ws.Range("J22").Value = 3500
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 3500
$ws.Range("M22").Value = -1805
$ws.Range("N22").Value = -4090

$ws.Range("H27").Value = 2800
$ws.Range("I27").Value = 2100
$ws.Range("J27").Value = 3500
$ws.Range("K27").Value = 2100
$ws.Range("L27").Value = 3500
$ws.Range("M27").Value = -1993
$ws.Range("N27").Value = -3714

$ws.Range("H40").Value = 2457
$ws.Range("J40").Value = 4165.6665
$ws.Range("L40").Value = 4165.6665
$ws.Range("N40").Value = -4437.6665

$ws.Range("H46").Value = 6400
$ws.Range("J46").Value = 6666.6665
$ws.Range("L46").Value = 6666.6665
$ws.Range("N46").Value = -7042.6665

$ws.Range("H55").Value = 1136.7333
$ws.Range("I55").Value = 1176
$ws.Range("J55").Value = 1077.8334
$ws.Range("K55").Value = 1176
$ws.Range("L55").Value = 1077.8334
$ws.Range("M55").Value = -1003
$ws.Range("N55").Value = -1423.8334

$ws.Range("H61").Value = 3986.1428
$ws.Range("I61").Value = 975.75
$ws.Range("K61").Value = 975.75
$ws.Range("M61").Value = -773.75

$ws.Range("H82").Value = 5834.375
$ws.Range("I82").Value = 4306.6665
$ws.Range("J82").Value = 6751
$ws.Range("K82").Value = 4306.6665
$ws.Range("L82").Value = 6751
$ws.Range("M82").Value = -3945.6665
$ws.Range("N82").Value = -7473

$ws.Range("H85").Value = 5834.375
$ws.Range("I85").Value = 4306.6665
$ws.Range("J85").Value = 6751
$ws.Range("K85").Value = 4306.6665
$ws.Range("L85").Value = 6751
$ws.Range("M85").Value = -3058.6665
$ws.Range("N85").Value = -9247

$ws.Range("H113").Value = 3986.1428
$ws.Range("I113").Value = 975.75
$ws.Range("K113").Value = 975.75
$ws.Range("M113").Value = 1194.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 19999.5
$ws.Range("J140").Value = 19999.5
$ws.Range("L140").Value = 19999.5
$ws.Range("N140").Value = -30359.5
